$p = $ppt.ActivePresentation
$p.Slides.Item(21).Delete()
